# Re-sort the song data table (B2:H42) by the Artist column (C) instead of
# the SongName column (B), keeping the ID column (A) fixed in place, and
# update the active cell selection to J31 to match the author's final
# interaction with the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("B2:H42")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C42"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("J31").Select()
